$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.44"
$ws.Range("D3").Value = "'28.19"
$ws.Range("E3").Value = "'-3.99%"
$ws.Range("D4").Value = "'5.270"
$ws.Range("E4").Value = "'1.53%"
$ws.Range("D5").Value = "'0.05707"
$ws.Range("E5").Value = "'-0.49%"
$ws.Range("D6").Value = "'6.640"
$ws.Range("E6").Value = "'1.30%"
$ws.Range("E7").Value = "'3.65%"
$ws.Range("D8").Value = "'0.8507"
$ws.Range("E8").Value = "'-0.85%"
$ws.Range("D9").Value = "'0.8825"
$ws.Range("E9").Value = "'1.56%"
$ws.Range("D10").Value = "'0.1387"
$ws.Range("E10").Value = "'1.67%"
$ws.Range("D11").Value = "'0.07078"
$ws.Range("E11").Value = "'-0.11%"
$ws.Range("D12").Value = "'0.03141"
$ws.Range("E12").Value = "'5.39%"
$ws.Range("D13").Value = "'0.09216"
$ws.Range("E13").Value = "'-1.80%"
$ws.Range("D14").Value = "'0.001530"
$ws.Range("E14").Value = "'-0.26%"
$ws.Range("D15").Value = "'0.0005957"
$ws.Range("E15").Value = "'-0.69%"
$ws.Range("D16").Value = "'0.005891"
$ws.Range("E16").Value = "'-2.17%"
$ws.Range("D17").Value = "'3.493"
$ws.Range("E17").Value = "'0.09%"
$ws.Range("E18").Value = "'0.20%"
$ws.Range("D19").Value = "'0.3168"
$ws.Range("E19").Value = "'-0.53%"
$ws.Range("D20").Value = "'0.03309"
$ws.Range("E20").Value = "'-2.58%"
$ws.Range("D21").Value = "'0.1306"
$ws.Range("E21").Value = "'1.59%"
$ws.Range("D22").Value = "'3.529"
$ws.Range("E22").Value = "'1.76%"
$ws.Range("D23").Value = "'0.04079"
$ws.Range("E23").Value = "'-1.32%"
$ws.Range("D24").Value = "'0.1378"
$ws.Range("E24").Value = "'-0.09%"
$ws.Range("D25").Value = "'0.001222"
$ws.Range("E25").Value = "'-0.21%"
$ws.Range("D26").Value = "'0.004154"
$ws.Range("E26").Value = "'-17.05%"
$ws.Range("D27").Value = "'0.0001199"
$ws.Range("E27").Value = "'-0.83%"
$ws.Range("D28").Value = "'0.0001448"
$ws.Range("D40").Value = "'0.03793"
$ws.Range("E40").Value = "'1.07%"
$ws.Range("D41").Value = "'0.1067"
$ws.Range("E41").Value = "'-0.62%"
$ws.Range("E42").Value = "'7.33%"
$ws.Range("E43").Value = "'-9.42%"
$ws.Range("D44").Value = "'0.009477"
$ws.Range("E44").Value = "'11.66%"
$ws.Range("D45").Value = "'0.00005280"
$ws.Range("E45").Value = "'0.40%"
$ws.Range("E46").Value = "'-0.03%"
$ws.Range("D47").Value = "'0.1049"
$ws.Range("E47").Value = "'84.15%"
$ws.Range("E48").Value = "'-0.33%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.03%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.03%"
